$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer / title cell (A1) with new timestamp
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 12:39"

# Row 14: Alemania
$ws.Range("A14").Value = "Alemania"
$ws.Range("B14").Value = 191216
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 174900
$ws.Range("E14").Value = 7355
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 8961

# Row 40: Suiza
$ws.Range("A40").Value = "Suiza"
$ws.Range("B40").Value = 31292
$ws.Range("C40").Value = 49
$ws.Range("D40").Value = 28900
$ws.Range("E40").Value = 436
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 1956

# Row 41: Filipinas
$ws.Range("A41").Value = "Filipinas"
$ws.Range("B41").Value = 30052
$ws.Range("C41").Value = 652
$ws.Range("D41").Value = 7893
$ws.Range("E41").Value = 20990
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 19
$ws.Range("H41").Value = 1169

# Row 42: Oman
$ws.Range("A42").Value = "Oman"
$ws.Range("B42").Value = 29471
$ws.Range("C42").Value = 905
$ws.Range("D42").Value = 15552
$ws.Range("E42").Value = 13788
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 131

# Row 48: Rumania
$ws.Range("A48").Value = "Rumania"
$ws.Range("B48").Value = 24045
$ws.Range("C48").Value = 315
$ws.Range("D48").Value = 16911
$ws.Range("E48").Value = 5622
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 12
$ws.Range("H48").Value = 1512

# Row 68: Marruecos
$ws.Range("A68").Value = "Marruecos"
$ws.Range("B68").Value = 9957
$ws.Range("C68").Value = 118
$ws.Range("D68").Value = 8249
$ws.Range("E68").Value = 1495
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 213

# Row 77: Consejo Danes para los Refugiados
$ws.Range("A77").Value = "Consejo Danes para los Refugiados"
$ws.Range("B77").Value = 5826
$ws.Range("C77").Value = 154
$ws.Range("D77").Value = 841
$ws.Range("E77").Value = 4855
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = 130

# Row 78: Senegal
$ws.Range("A78").Value = "Senegal"
$ws.Range("B78").Value = 5783
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 3859
$ws.Range("E78").Value = 1842
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 82

# Row 83: El Salvador
$ws.Range("A83").Value = "El Salvador"
$ws.Range("B83").Value = 4626
$ws.Range("C83").Value = 151
$ws.Range("D83").Value = 2535
$ws.Range("E83").Value = 1993
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = 98

# Row 84: Republica de Yibuti
$ws.Range("A84").Value = "Republica de Yibuti"
$ws.Range("B84").Value = 4565
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 3565
$ws.Range("E84").Value = 955
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 45

# Row 85: Kenia
$ws.Range("A85").Value = "Kenia"
$ws.Range("B85").Value = 4478
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 1586
$ws.Range("E85").Value = 2771
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 121

# Row 106: Albania
$ws.Range("A106").Value = "Albania"
$ws.Range("B106").Value = 1962
$ws.Range("C106").Value = 71
$ws.Range("D106").Value = 1134
$ws.Range("E106").Value = 784
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 44

# Row 107: Sri Lanka
$ws.Range("A107").Value = "Sri Lanka"
$ws.Range("B107").Value = 1950
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 1498
$ws.Range("E107").Value = 441
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 11

# Row 108: Mali
$ws.Range("A108").Value = "Mali"
$ws.Range("B108").Value = 1933
$ws.Range("C108").Value = 10
$ws.Range("D108").Value = 1255
$ws.Range("E108").Value = 569
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 109

# Row 114: Madagascar
$ws.Range("A114").Value = "Madagascar"
$ws.Range("B114").Value = 1596
$ws.Range("C114").Value = 93
$ws.Range("D114").Value = 655
$ws.Range("E114").Value = 927
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 14

# Row 115: Eslovaquia
$ws.Range("A115").Value = "Eslovaquia"
$ws.Range("B115").Value = 1587
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 1447
$ws.Range("E115").Value = 112
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 28

# Row 116: Guinea-Bisau
$ws.Range("A116").Value = "Guinea-Bisau"
$ws.Range("B116").Value = 1541
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 153
$ws.Range("E116").Value = 1371
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 17

# Row 117: Libano
$ws.Range("A117").Value = "Libano"
$ws.Range("B117").Value = 1536
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 1006
$ws.Range("E117").Value = 498
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 32

# Row 118: Eslovenia
$ws.Range("A118").Value = "Eslovenia"
$ws.Range("B118").Value = 1525
$ws.Range("C118").Value = 6
$ws.Range("D118").Value = 1359
$ws.Range("E118").Value = 57
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 109

# Row 119: Nueva Zelanda
$ws.Range("A119").Value = "Nueva Zelanda"
$ws.Range("B119").Value = 1511
$ws.Range("C119").Value = 2
$ws.Range("D119").Value = 1482
$ws.Range("E119").Value = 7
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 22

# Row 124: Hong Kong
$ws.Range("A124").Value = "Hong Kong"
$ws.Range("B124").Value = 1132
$ws.Range("C124").Value = 3
$ws.Range("D124").Value = 1078
$ws.Range("E124").Value = 49
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 5

# Row 139: Malaui
$ws.Range("A139").Value = "Malaui"
$ws.Range("B139").Value = 730
$ws.Range("C139").Value = 110
$ws.Range("D139").Value = 91
$ws.Range("E139").Value = 628
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 3
$ws.Range("H139").Value = 11

# Row 140: Crucero
$ws.Range("A140").Value = "Crucero"
$ws.Range("B140").Value = 712
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 651
$ws.Range("E140").Value = 48
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 13

# Row 141: Ruanda
$ws.Range("A141").Value = "Ruanda"
$ws.Range("B141").Value = 702
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 357
$ws.Range("E141").Value = 343
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 2

# Row 142: Santo Tome y Principe
$ws.Range("A142").Value = "Santo Tome y Principe"
$ws.Range("B142").Value = 698
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 203
$ws.Range("E142").Value = 483
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 12

# Row 143: San Marino
$ws.Range("A143").Value = "San Marino"
$ws.Range("B143").Value = 696
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 610
$ws.Range("E143").Value = 44
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 42

# Row 144: Mozambique
$ws.Range("A144").Value = "Mozambique"
$ws.Range("B144").Value = 688
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 177
$ws.Range("E144").Value = 507
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 4

# Row 145: Malta
$ws.Range("A145").Value = "Malta"
$ws.Range("B145").Value = 664
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 616
$ws.Range("E145").Value = 39
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 9

# Row 146: Jamaica
$ws.Range("A146").Value = "Jamaica"
$ws.Range("B146").Value = 657
$ws.Range("C146").Value = 5
$ws.Range("D146").Value = 462
$ws.Range("E146").Value = 185
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 10

# Row 147: Benin
$ws.Range("A147").Value = "Benin"
$ws.Range("B147").Value = 650
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 247
$ws.Range("E147").Value = 392
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 11

# Row 148: Suazilandia
$ws.Range("A148").Value = "Suazilandia"
$ws.Range("B148").Value = 627
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 285
$ws.Range("E148").Value = 337
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 5

# Row 165: Siria
$ws.Range("A165").Value = "Siria"
$ws.Range("B165").Value = 204
$ws.Range("C165").Value = 6
$ws.Range("D165").Value = 83
$ws.Range("E165").Value = 114
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 7

# Row 173: Burundi
$ws.Range("A173").Value = "Burundi"
$ws.Range("B173").Value = 144
$ws.Range("C173").Value = 40
$ws.Range("D173").Value = 93
$ws.Range("E173").Value = 50
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 1

# Row 174: Eritrea
$ws.Range("A174").Value = "Eritrea"
$ws.Range("B174").Value = 143
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 39
$ws.Range("E174").Value = 104
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# Row 175: Brunei
$ws.Range("A175").Value = "Brunei"
$ws.Range("B175").Value = 141
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 138
$ws.Range("E175").Value = 0
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 3

# Row 176: Camboya
$ws.Range("A176").Value = "Camboya"
$ws.Range("B176").Value = 129
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 127
$ws.Range("E176").Value = 2
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

# Row 177: Trinidad yTobago
$ws.Range("A177").Value = "Trinidad yTobago"
$ws.Range("B177").Value = 123
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 109
$ws.Range("E177").Value = 6
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 8

# Row 202: Fiyi
$ws.Range("A202").Value = "Fiyi"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

# Row 203: Dominica
$ws.Range("A203").Value = "Dominica"
$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 208: Santa Sede
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 12
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# Row 209: Islas Turcas y Caicos
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("B209").Value = 12
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 1

# Row 213: Papua Nueva Guinea
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: Islas Virgenes Britanicas
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
